# Schedule_Info.xlsx update — "Add files via upload"
#
# Grade_1 (sheet1): message #4 (row 4) has now been sent -> Status=1,
#   Message_Id=7488; the three trailing placeholder rows (105-107) are removed.
# Grade_2 (sheet2): message #4 (row 4) sent -> Status=1, Message_Id=3461;
#   4 new upcoming rows appended (99-102).
# Grade_3 (sheet3): message #4 (row 4) sent -> Status=1, Message_Id=2494;
#   5 new upcoming rows appended (97-101).
# Grade_4 (sheet4): message #4 (row 4) sent -> Status=1, Message_Id=6226;
#   7 new upcoming rows appended (96-102).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Grade_1
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Grade_1")
$ws1.Cells.Item(4,3).Value() = 1
$ws1.Cells.Item(4,4).Value() = 7488
$ws1.Rows("105:107").Delete()

# ---------------------------------------------------------------------------
# Grade_2
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Grade_2")
$ws2.Cells.Item(4,3).Value() = 1
$ws2.Cells.Item(4,4).Value() = 3461

$grade2New = @(
    @("16,03,2023", "Grade2_(98).png"),
    @("17,03,2023", "Grade2_(99).png"),
    @("18,03,2023", "Grade2_(100).png"),
    @("19,03,2023", "Grade2_(101).png")
)
$r = 99
foreach ($entry in $grade2New) {
    $ws2.Cells.Item($r,1).Value() = $entry[0]
    $ws2.Cells.Item($r,2).Value() = $entry[1]
    $ws2.Cells.Item($r,3).Value() = 0
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Grade_3
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Grade_3")
$ws3.Cells.Item(4,3).Value() = 1
$ws3.Cells.Item(4,4).Value() = 2494

$grade3New = @(
    @("14,03,2023", "Grade3_(96).png"),
    @("15,03,2023", "Grade3_(97).png"),
    @("16,03,2023", "Grade3_(98).png"),
    @("17,03,2023", "Grade3_(99).png"),
    @("18,03,2023", "Grade3_(100).png")
)
$r = 97
foreach ($entry in $grade3New) {
    $ws3.Cells.Item($r,1).Value() = $entry[0]
    $ws3.Cells.Item($r,2).Value() = $entry[1]
    $ws3.Cells.Item($r,3).Value() = 0
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Grade_4
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Grade_4")
$ws4.Cells.Item(4,3).Value() = 1
$ws4.Cells.Item(4,4).Value() = 6226

$grade4New = @(
    @("13,03,2023", "Grade4_(95).png"),
    @("14,03,2023", "Grade4_(96).png"),
    @("15,03,2023", "Grade4_(97).png"),
    @("16,03,2023", "Grade4_(98).png"),
    @("17,03,2023", "Grade4_(99).png"),
    @("18,03,2023", "Grade4_(100).png"),
    @("19,03,2023", "Grade4_(101).png")
)
$r = 96
foreach ($entry in $grade4New) {
    $ws4.Cells.Item($r,1).Value() = $entry[0]
    $ws4.Cells.Item($r,2).Value() = $entry[1]
    $ws4.Cells.Item($r,3).Value() = 0
    $r = $r + 1
}
